$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$newSheet = $wb.Worksheets.Add($null, $ws1)
$newSheet.Name = "Emp Utilization"
$newSheet.Range("A4").Value = "Employees"
[void]$newSheet.Range("A4").Select()
[void]$ws1.Activate()
